$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = "https://dev.i2-donate.com/i2D-Publish-Docs/i2-Donate Terms and Conditions.html"
$ws.Range("B4").Value = "https://dev.i2-donate.com/i2D-Publish-Docs/i2-Donate Privacy Policy.html"
$ws.Range("B5").Value = "https://dev.i2-donate.com/i2D-Publish-Docs/i2-Donate Help and Support.html"
$ws.Range("B6").Value = "https://dev.i2-donate.com/i2D-Publish-Docs/i2-Donate About App.html"

$ws.Range("B6").Select()
